$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (row 2 and row 3) represent the same market entry but for
# two different weekly dates. The weekly update swaps which week's figures sit
# on top: row 2 now holds the figures that used to be in row 3, and vice versa.

# Row 2 gets the values that were previously in row 3
$ws.Range("D2").Value = 44875
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1600
$ws.Range("L2").Value = 1700
$ws.Range("M2").Value = 1650
$ws.Range("P2").Value = 1650

# Row 3 gets the values that were previously in row 2
$ws.Range("D3").Value = 44547
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1600
$ws.Range("M3").Value = 1550
$ws.Range("P3").Value = 1550
